$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 12 : id 11 ----
$ws.Cells.Item(12, 1).Value = 11
$ws.Cells.Item(12, 2).Value = 44231
$ws.Cells.Item(11, 2).Copy()
$ws.Cells.Item(12, 2).PasteSpecial(-4122)
$ws.Cells.Item(12, 3).Value = "Google Scholar"
$ws.Cells.Item(12, 4).Value = "lemon market asymmetric information"
$ws.Cells.Item(12, 5).Value = "https://www.tandfonline.com/doi/pdf/10.1300/J130v05n01_02?needAccess=true"

# ---- Row 13 : id 12 (link entered before keyword, matching original authoring order) ----
$ws.Cells.Item(13, 1).Value = 12
$ws.Cells.Item(13, 2).Value = 44231
$ws.Cells.Item(11, 2).Copy()
$ws.Cells.Item(13, 2).PasteSpecial(-4122)
$ws.Cells.Item(13, 3).Value = "Google Scholar"
$ws.Cells.Item(13, 5).Value = "https://www.jstor.org/stable/1810022?casa_token=fVot-7ZHzHcAAAAA%3Av2RSbb-ydlHNPCWKqzpiBi9xDE178SO7x0r4brwKhlxBz2zwLNIGy6Bt1_fgrVdHRIFyyCL2CuirNEXWp7DOod2eiFDAqER6BDBPUiKnEFRxzA9vcPNR&seq=1#metadata_info_tab_contents"
$ws.Cells.Item(13, 4).Value = "lemon car market"

# ---- Row 14 : id 13 ----
$ws.Cells.Item(14, 1).Value = 13
$ws.Cells.Item(14, 2).Value = 44231
$ws.Cells.Item(11, 2).Copy()
$ws.Cells.Item(14, 2).PasteSpecial(-4122)
$ws.Cells.Item(14, 3).Value = "Google Scholar"
$ws.Cells.Item(14, 4).Value = "George akerlof the market for lemons"
$ws.Cells.Item(14, 5).Value = "https://www.williamdavid.me.uk/wp-content/uploads/2017/04/The-Market-for-Lemons-text-full.pdf"

# ---- Row 15 : id 14 ----
$ws.Cells.Item(15, 1).Value = 14
$ws.Cells.Item(15, 2).Value = 44231
$ws.Cells.Item(11, 2).Copy()
$ws.Cells.Item(15, 2).PasteSpecial(-4122)
$ws.Cells.Item(15, 3).Value = "Google Scholar"
$ws.Cells.Item(15, 4).Value = "George akerlof the market for lemons"
$ws.Cells.Item(15, 5).Value = "https://www.jstor.org/stable/pdf/1879431.pdf?refreqid=excelsior%3A1ff10c03269b5441b7568b90db57b572"

$ws.Range("E7").Select()
